# Refresh the "cryptos" price/volume snapshot (GitHub Actions bot update).
# Column D ("Price") and column E ("Volume(1h)") are stored as plain text in
# this sheet (e.g. "308.28", "  +0.05%  "), so cells whose new value would
# otherwise be auto-parsed as a number by Excel are written with a leading
# apostrophe to keep them as text, matching the original cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.270.40"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").Value = "2.271.05"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'308.76"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("D6").Value = "'97.13"
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("E7").Value = "  -1.11%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -1.54%  "
$ws.Range("D10").Value = "'35.15"
$ws.Range("E10").Value = "  -3.46%  "
$ws.Range("D11").Value = "'0.0808"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("D13").Value = "'6.79"
$ws.Range("E13").Value = "  +0.68%  "
$ws.Range("D14").Value = "2.622.50"
$ws.Range("E14").Value = "  -1.11%  "
$ws.Range("D15").Value = "'14.64"
$ws.Range("E15").Value = "  -0.02%  "
$ws.Range("D16").Value = "2.269.76"
$ws.Range("E16").Value = "  -1.46%  "
$ws.Range("E17").Value = "  -1.86%  "
$ws.Range("D18").Value = "42.113.63"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("D19").Value = "'12.25"
$ws.Range("E19").Value = "  -3.94%  "
$ws.Range("D20").Value = "0.0₃0906"
$ws.Range("E20").Value = "  -1.75%  "
$ws.Range("E21").Value = "  -1.09%  "
$ws.Range("D22").Value = "'67.70"
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("E23").Value = "  -2.91%  "
$ws.Range("D24").Value = "'2.59"
$ws.Range("E24").Value = "  -0.85%  "
$ws.Range("E25").Value = "  +0.69%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").Value = "'23.55"
$ws.Range("E27").Value = "  -2.19%  "
$ws.Range("D28").Value = "'36.66"
$ws.Range("E28").Value = "  -0.62%  "
$ws.Range("E29").Value = "  -0.39%  "
$ws.Range("E30").Value = "  +0.58%  "
$ws.Range("E31").Value = "  +2.09%  "
$ws.Range("E32").Value = "  -1.34%  "
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("E34").Value = "  -0.60%  "
$ws.Range("D35").Value = "'0.0736"
$ws.Range("E35").Value = "  -2.27%  "
$ws.Range("D36").Value = "'17.51"
$ws.Range("E36").Value = "  +0.65%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("E38").Value = "  -3.70%  "
$ws.Range("E39").Value = "  -3.72%  "
$ws.Range("E40").Value = "  -0.97%  "
$ws.Range("D41").Value = "'4.14"
$ws.Range("E41").Value = "  -1.32%  "
$ws.Range("D42").Value = "'2.28"
$ws.Range("E42").Value = "  -6.33%  "
$ws.Range("D43").Value = "1.953.91"
$ws.Range("E43").Value = "  -2.66%  "
$ws.Range("D44").Value = "'18.93"
$ws.Range("E44").Value = "  -1.56%  "
$ws.Range("E45").Value = "  -1.93%  "
$ws.Range("E46").Value = "  -3.43%  "
$ws.Range("D47").Value = "'9.82"
$ws.Range("E47").Value = "  -4.69%  "
$ws.Range("D48").Value = "'53.59"
$ws.Range("E48").Value = "  -0.78%  "
$ws.Range("D49").Value = "2.494.29"
$ws.Range("E49").Value = "  -1.03%  "
$ws.Range("D50").Value = "'92.38"
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("D51").Value = "'71.56"
$ws.Range("E51").Value = "  -1.70%  "
